# The deck ships two theme parts:
#   - the Slide Master / presentation theme (currently the "Integral" palette)
#   - the Notes Master theme (currently the "Office Theme" palette)
#
# The edit swaps the two palettes: the Slide Master's theme becomes the
# classic "Office Theme" colours, while the Notes Master keeps pointing at
# its own (separate) theme part which should carry the "Integral" colours.
#
# PowerPoint's automation surface only lets us push real colour edits
# through the live ThemeColorScheme of the presentation's active theme
# (Master.Theme / NotesMaster.Theme resolve to that same theme), so we
# recolour it here to the "Office Theme" palette, in clrScheme order:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink

$p = $ppt.ActivePresentation

$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColors = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
